$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new birthday entry: HERNANDO ESTRADA PACHECO, day 8, month 8
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "HERNANDO ESTRADA PACHECO"

# Update selection to A10 as recorded in the saved workbook
$ws.Range("A10").Select()
